$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update the test-case name cells (PR_ -> PM_ prefix rename)
$ws.Range("B8").Value = "TC_PM_COA_SEC_AgencyGLAccountProfile_AddNewProfile"
$ws.Range("B10").Value = "TC_PM_COA_SEC_AgencyGLAccountProfile_EditProfile"
$ws.Range("B13").Value = "TC_PM_COA_SEC_AgencyGLAccountProfile_AddNewProfile"

# Adjust column widths for B and C (bestFit-style autofit after the rename)
$ws.Columns.Item(2).ColumnWidth = 49.833333333333333
$ws.Columns.Item(3).ColumnWidth = 39.5

# Update selection to C18
$ws.Range("C18").Select()
